$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.183.13"
$ws.Range("E2").Value = "  -1.90%  "

$ws.Range("D3").Value = "3.074.40"
$ws.Range("E3").Value = "  -1.61%  "

$ws.Range("E4").Value = "  +0.00%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "522.08"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.34%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "135.59"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -4.92%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "3.074.53"
$ws.Range("E8").Value = "  -1.58%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.468"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +5.01%  "

$ws.Range("E10").Value = "  +1.58%  "

$ws.Range("E11").Value = "  -2.86%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.400"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.60%  "

$ws.Range("E13").Value = "  +1.04%  "

$ws.Range("D14").Value = "3.595.79"
$ws.Range("E14").Value = "  -1.75%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "25.16"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.98%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000160"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -3.29%  "

$ws.Range("D17").Value = "57.257.37"
$ws.Range("E17").Value = "  -1.82%  "

$ws.Range("D18").Value = "3.075.35"
$ws.Range("E18").Value = "  -1.51%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "5.86"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -4.26%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.40"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.49%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.81"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.24%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "348.97"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.63%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "69.01"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.06%  "

$ws.Range("E25").Value = "  -3.19%  "

$ws.Range("E26").Value = "  -2.68%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.27%  "

$ws.Range("D28").Value = "0.0₃0858"
$ws.Range("E28").Value = "  -8.36%  "

$ws.Range("E29").Value = "  -0.02%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.17"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.70%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.85"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.13%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "20.89"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.12%  "

$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.76"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -10.13%  "

$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "159.31"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.52%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.83"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.48%  "

$ws.Range("E36").Value = "  -4.83%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.99"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.04%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "25.34"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.62%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.23"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.63%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0655"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.85%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "4.03"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.15%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.55"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -6.79%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.692"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("D44").Value = "2.410.59"
$ws.Range("E44").Value = "  +5.83%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "36.48"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.48%  "

$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").Value = "3.112.61"
$ws.Range("E47").Value = "  -1.61%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0261"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.55%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "5.95"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.63%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.934"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -6.89%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "19.53"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -5.43%  "
